$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$acid = "ACID adalah sekumpulan properti (Atomicity, Consistency, Isolation, Durability) yang menjamin keandalan transaksi dalam basis data."
$ai = "AI adalah Artificial Intelligence, yaitu kecerdasan buatan yang memungkinkan mesin untuk meniru kecerdasan manusia dalam menyelesaikan tugas."

# New question/answer pairs for rows 2-15
$questions = @(
    "api?",
    "oop",
    "asiap",
    "commit dalam git?",
    "bahasa pemrograman tingkat rendah",
    "debungging",
    "apa itu go",
    "ai",
    "ai",
    "ai",
    "ai",
    "acid",
    "debungging",
    "commit dalam git?"
)

$answers = @(
    $acid,
    $acid,
    $acid,
    $ai,
    $ai,
    $acid,
    $acid,
    $acid,
    $acid,
    $acid,
    $acid,
    $ai,
    $acid,
    $ai
)

for ($i = 0; $i -lt $questions.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $questions[$i]
    $ws.Cells.Item($row, 2).Value = $answers[$i]
}

# Remove rows 16 and 17 (previously siem/metadata), shrinking the range to A1:B15
$ws.Range("A16:B17").Delete()

$wb.Save()
